# The vendor "UPS" row (row 6) had a stray "x" mark in column B (the
# "use this vendor" checkbox column) and a leftover Bill Amount of 12.43
# in column D that were left over from testing and were messing up the
# printed invoice. Clear both so the row matches the other vendor rows
# (which only have the vendor name in column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").ClearContents()
$ws.Range("D6").ClearContents()

# Leave the selection where the user ended up after cleaning the row.
$ws.Range("B6").Select()
